$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-converted to numbers by Excel (losing trailing zeros / exact text).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.693.10"
$ws.Range("E2").Value = "  -2.07%  "
$ws.Range("D3").Value = "1.754.70"
$ws.Range("E3").Value = "  -2.71%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "324.08"
$ws.Range("E5").Value = "  -4.56%  "
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "0.4299"
$ws.Range("E7").Value = "  -7.30%  "
$ws.Range("D8").Value = "0.3640"
$ws.Range("E8").Value = "  -4.50%  "
$ws.Range("D9").Value = "45.31"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "0.07480"
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("E11").Value = "  -3.33%  "
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "21.59"
$ws.Range("E13").Value = "  -4.11%  "
$ws.Range("E14").Value = "  -3.30%  "
$ws.Range("D15").Value = "7.246"
$ws.Range("E15").Value = "  -4.19%  "
$ws.Range("D16").Value = "1.751.20"
$ws.Range("E16").Value = "  -3.14%  "
$ws.Range("E17").Value = "  -2.39%  "
$ws.Range("D18").Value = "87.97"
$ws.Range("E18").Value = "  +7.79%  "
$ws.Range("D19").Value = "0.06214"
$ws.Range("E19").Value = "  -7.55%  "
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").Value = "17.11"
$ws.Range("E21").Value = "  -2.36%  "
$ws.Range("D22").Value = "6.149"
$ws.Range("E22").Value = "  -4.54%  "
$ws.Range("D23").Value = "0.5273"
$ws.Range("E23").Value = "  -4.62%  "
$ws.Range("D24").Value = "27.718.69"
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("D25").Value = "11.65"
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("D26").Value = "2.324"
$ws.Range("E26").Value = "  -4.10%  "
$ws.Range("D27").Value = "20.54"
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("D28").Value = "152.55"
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("D29").Value = "2.365"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").Value = "1.950.46"
$ws.Range("E30").Value = "  -3.14%  "
$ws.Range("D31").Value = "1.220"
$ws.Range("E31").Value = "  -2.72%  "
$ws.Range("D32").Value = "127.30"
$ws.Range("E32").Value = "  -4.35%  "
$ws.Range("D33").Value = "5.718"
$ws.Range("E33").Value = "  -2.33%  "
$ws.Range("D34").Value = "0.09146"
$ws.Range("E34").Value = "  -4.51%  "
$ws.Range("D35").Value = "3.657"
$ws.Range("E35").Value = "  -9.47%  "
$ws.Range("D36").Value = "12.70"
$ws.Range("E36").Value = "  +4.82%  "
$ws.Range("D37").Value = "0.02309"
$ws.Range("E37").Value = "  -1.96%  "
$ws.Range("D38").Value = "0.2155"
$ws.Range("E38").Value = "  -6.57%  "
$ws.Range("D39").Value = "5.109"
$ws.Range("E39").Value = "  -3.28%  "
$ws.Range("D40").Value = "0.6479"
$ws.Range("E40").Value = "  -2.37%  "
$ws.Range("D41").Value = "0.06098"
$ws.Range("E41").Value = "  -4.05%  "
$ws.Range("D42").Value = "1.197"
$ws.Range("E42").Value = "  -3.70%  "
$ws.Range("D43").Value = "1.427"
$ws.Range("E43").Value = "  -4.18%  "
$ws.Range("D44").Value = "7.951"
$ws.Range("E44").Value = "  -5.18%  "
$ws.Range("D45").Value = "0.9988"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("D46").Value = "13.72"
$ws.Range("E46").Value = "  -3.08%  "
$ws.Range("D47").Value = "3.745"
$ws.Range("E47").Value = "  -3.37%  "
$ws.Range("D48").Value = "0.5937"
$ws.Range("E48").Value = "  -3.38%  "
$ws.Range("D49").Value = "125.97"
$ws.Range("E49").Value = "  -3.84%  "
$ws.Range("D50").Value = "1.976"
$ws.Range("E50").Value = "  -3.18%  "
$ws.Range("D51").Value = "0.06899"
$ws.Range("E51").Value = "  -4.19%  "
